# Add a new "performance_glitch_user" test row to the DDT testdata sheet.
#
# Row 4:
#   A4 = "performance_glitch_user" with a new highlighted style
#        (white "DM Mono" 11pt text on a dark teal/green fill)
#   B4 = "secret_sauce" using the existing default data-row style (same as B3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A4: new value + new font/fill formatting ---
$ws.Range("A4").Value = "performance_glitch_user"
$ws.Range("A4").Font.Name = '"DM Mono"'
$ws.Range("A4").Font.Size = 11
$ws.Range("A4").Font.Color = 16777215       # white (RGB 255,255,255 -> 0xFFFFFF)
$ws.Range("A4").Interior.Color = 2237203    # dark teal fill (RGB 0x13,0x23,0x22 -> BGR 0x222313)

# --- B4: new value, reuse the same style already used by B2/B3 (secret_sauce) ---
# Copying a formatted cell and then overwriting its value keeps the destination
# on the existing shared cell style instead of minting a new one.
$ws.Range("B3").Copy($ws.Range("B4"))
$ws.Range("B4").Value = "secret_sauce"
